$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '24.744.97'
$ws.Range("E2").Value = '  +1.78%  '

# Row 3
$ws.Range("D3").Value = '1.662.86'
$ws.Range("E3").Value = '  +1.94%  '

# Row 4
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").Value = '329.78'

# Row 6
$ws.Range("E6").Value = '  -0.17%  '

# Row 7
$ws.Range("D7").Value = '0.3651'
$ws.Range("E7").Value = '  +1.30%  '

# Row 8
$ws.Range("D8").Value = '47.35'
$ws.Range("E8").Value = '  +1.21%  '

# Row 9
$ws.Range("D9").Value = '0.3237'
$ws.Range("E9").Value = '  +0.40%  '

# Row 10
$ws.Range("D10").Value = '1.134'
$ws.Range("E10").Value = '  +3.26%  '

# Row 11
$ws.Range("D11").Value = '0.07050'
$ws.Range("E11").Value = '  +3.04%  '

# Row 12
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  -0.07%  '

# Row 13
$ws.Range("D13").Value = '6.056'
$ws.Range("E13").Value = '  +3.29%  '

# Row 14
$ws.Range("D14").Value = '19.53'
$ws.Range("E14").Value = '  +3.29%  '

# Row 15
$ws.Range("D15").Value = '1.660.87'
$ws.Range("E15").Value = '  +1.58%  '

# Row 16
$ws.Range("D16").Value = '6.589'
$ws.Range("E16").Value = '  +1.88%  '

# Row 17
$ws.Range("D17").Value = '0.00001049'
$ws.Range("E17").Value = '  +1.49%  '

# Row 18
$ws.Range("D18").Value = '0.06628'
$ws.Range("E18").Value = '  +2.27%  '

# Row 19
$ws.Range("E19").Value = '  -0.14%  '

# Row 20
$ws.Range("D20").Value = '78.51'
$ws.Range("E20").Value = '  +3.45%  '

# Row 21
$ws.Range("D21").Value = '5.930'
$ws.Range("E21").Value = '  +1.57%  '

# Row 22
$ws.Range("D22").Value = '15.79'
$ws.Range("E22").Value = '  +1.86%  '

# Row 23
$ws.Range("D23").Value = '12.46'
$ws.Range("E23").Value = '  +4.38%  '

# Row 24
$ws.Range("D24").Value = '24.765.88'
$ws.Range("E24").Value = '  +2.05%  '

# Row 25
$ws.Range("D25").Value = '2.457'
$ws.Range("E25").Value = '  +2.41%  '

# Row 26
$ws.Range("D26").Value = '2.413'
$ws.Range("E26").Value = '  +6.36%  '

# Row 27
$ws.Range("D27").Value = '148.48'
$ws.Range("E27").Value = '  +2.69%  '

# Row 28
$ws.Range("E28").Value = '  +1.70%  '

# Row 29
$ws.Range("D29").Value = '1.845.99'
$ws.Range("E29").Value = '  +1.60%  '

# Row 30
$ws.Range("D30").Value = '125.79'
$ws.Range("E30").Value = '  +2.56%  '

# Row 31
$ws.Range("D31").Value = '1.163'
$ws.Range("E31").Value = '  +2.75%  '

# Row 32
$ws.Range("D32").Value = '4.070'
$ws.Range("E32").Value = '  +0.58%  '

# Row 33
$ws.Range("D33").Value = '5.693'
$ws.Range("E33").Value = '  +3.81%  '

# Row 34
$ws.Range("D34").Value = '0.08499'
$ws.Range("E34").Value = '  +2.41%  '

# Row 35
$ws.Range("D35").Value = '1.637'
$ws.Range("E35").Value = '  -1.08%  '

# Row 36
$ws.Range("D36").Value = '12.11'
$ws.Range("E36").Value = '  -0.78%  '

# Row 37
$ws.Range("D37").Value = '5.146'
$ws.Range("E37").Value = '  +1.57%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02271'
$ws.Range("E38").Value = '  +4.02%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.06104'
$ws.Range("E39").Value = '  +2.57%  '

# Row 40
$ws.Range("D40").Value = '1.222'
$ws.Range("E40").Value = '  +2.41%  '

# Row 41
$ws.Range("D41").Value = '0.2080'
$ws.Range("E41").Value = '  +3.52%  '

# Row 42
$ws.Range("D42").Value = '8.198'

# Row 43
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  -0.10%  '

# Row 44
$ws.Range("D44").Value = '0.5909'
$ws.Range("E44").Value = '  +2.37%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '13.48'
$ws.Range("E45").Value = '  +8.13%  '

# Row 46
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = '3.848'
$ws.Range("E46").Value = '  +3.86%  '

# Row 47
$ws.Range("D47").Value = '0.5650'
$ws.Range("E47").Value = '  +3.03%  '

# Row 48
$ws.Range("D48").Value = '125.37'
$ws.Range("E48").Value = '  +4.12%  '

# Row 49
$ws.Range("D49").Value = '1.944'
$ws.Range("E49").Value = '  +1.97%  '

# Row 50
$ws.Range("D50").Value = '0.06962'
$ws.Range("E50").Value = '  +1.68%  '

# Row 51
$ws.Range("E51").Value = '  +5.63%  '
